$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Intro paragraph: describe labels for the two input fields and
#    relocate the hidden "_GoBack" bookmark to sit right after the
#    new "Player 2" default-value sentence.
# ---------------------------------------------------------------

$find = $d.Content.Find
$find.Execute(
    " game’s index page where are two input fields and a button. Input field are for inserting player names. First input field has by default “Player 1”, second input field “Player 2”. Button has label “Start game”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " game’s index page where are two input fields with respective labels and a button. Input field are for inserting player names. First input field has a label “Player 1” and the default value “Player 1” as a name as well, second input field has a label “Player 2” and default value “Player 2”. Button has label “Start game”.",
    2) | Out-Null

# Move the hidden "_GoBack" bookmark here (right before ". Button has label").
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$find2 = $d.Content.Find
$find2.Execute("default value “Player 2”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$gobackPos = $find2.Parent.End
$d.Bookmarks.Add("_GoBack", $d.Range($gobackPos, $gobackPos)) | Out-Null

# ---------------------------------------------------------------
# 2) "Undo" passage: the lastRenderedPageBreak marker shifts earlier
#    in the sentence (cosmetic repagination). Re-word the three
#    sentences around it without crossing the run that carries the
#    page-break marker, so the marker naturally ends up attached to
#    the correct text.
# ---------------------------------------------------------------

$find3 = $d.Content.Find
$find3.Execute("Game moved last moved checker from point 10 back to point 13. ", $true, $false, $false, $false, $false, $true, 1, $false, "Game moved last moved ", 2) | Out-Null

$find4 = $d.Content.Find
$find4.Execute("Mart pushed again “Undo” and game moved ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

$find5 = $d.Content.Find
$find5.Execute("the firstly moved checker from point 10 to point 13.", $true, $false, $false, $false, $false, $true, 1, $false, "checker from point 10 back to point 13. Mart pushed again “Undo” and game moved the firstly moved checker from point 10 to point 13.", 2) | Out-Null

# ---------------------------------------------------------------
# 3) "Congratulation" passage: the lastRenderedPageBreak marker
#    shifts earlier as well.
# ---------------------------------------------------------------

$find6 = $d.Content.Find
$find6.Execute("Since it is the last checker, Mart won the game. ", $true, $false, $false, $false, $false, $true, 1, $false, "Since it is ", 2) | Out-Null

$find7 = $d.Content.Find
$find7.Execute("On the board pops up a div with label ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

$find8 = $d.Content.Find
$find8.Execute("“Congratulation Mart you won!” and with buttons “", $true, $false, $false, $false, $false, $true, 1, $false, "the last checker, Mart won the game. On the board pops up a div with label “Congratulation Mart you won!” and with buttons “", 2) | Out-Null
